$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.106.44'
$ws.Range("E2").Value = '  -2.30%  '
$ws.Range("D3").Value = '2.136.44'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.18'
$ws.Range("E5").Value = '  -3.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.597'
$ws.Range("E6").Value = '  -4.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.73'
$ws.Range("E7").Value = '  -6.36%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -7.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.02'
$ws.Range("E10").Value = '  -11.29%  '
$ws.Range("E11").Value = '  -7.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.10'
$ws.Range("E12").Value = '  -7.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0991'
$ws.Range("E13").Value = '  -4.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.53'
$ws.Range("E14").Value = '  -7.57%  '
$ws.Range("D15").Value = '2.457.88'
$ws.Range("E15").Value = '  -3.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '2.126.37'
$ws.Range("E17").Value = '  -3.65%  '
$ws.Range("E18").Value = '  -8.06%  '
$ws.Range("D19").Value = '40.997.50'
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("E20").Value = '  -8.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.78'
$ws.Range("E21").Value = '  -5.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.67'
$ws.Range("E22").Value = '  -8.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '223.12'
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.42'
$ws.Range("E24").Value = '  -14.28%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.88'
$ws.Range("E26").Value = '  -10.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.46'
$ws.Range("E27").Value = '  -11.58%  '
$ws.Range("E28").Value = '  -9.02%  '
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("E30").Value = '  -6.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.75'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("E32").Value = '  -5.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.47'
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0743'
$ws.Range("E34").Value = '  -7.03%  '
$ws.Range("E35").Value = '  -12.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0992'
$ws.Range("E37").Value = '  -10.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.05'
$ws.Range("E38").Value = '  -5.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0279'
$ws.Range("E39").Value = '  -8.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  -5.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.49'
$ws.Range("E41").Value = '  -19.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.20'
$ws.Range("E42").Value = '  -8.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '56.66'
$ws.Range("E43").Value = '  -13.38%  '
$ws.Range("E44").Value = '  -7.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.12'
$ws.Range("E45").Value = '  -8.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0948'
$ws.Range("E46").Value = '  -5.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.99'
$ws.Range("E47").Value = '  -8.78%  '
$ws.Range("E48").Value = '  -5.24%  '
$ws.Range("E49").Value = '  -6.72%  '
$ws.Range("E50").Value = '  -3.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.11'
$ws.Range("E51").Value = '  -12.57%  '
